$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the column headers to use the respective input-file-name suffix
#    (_old -> _FV2304, _new -> _FV2310) instead of the generic _old/_new.
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Text -replace '_old$', '_FV2304')
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Text -replace '_new$', '_FV2310')
}

# 2. Turn the used range into a real table ("Table1") with an AutoFilter.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# 3. Freeze the header row (split below row 1) and select the top-left cell
#    of the scrollable (bottom-left) pane, like the authored workbook does.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
